$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.889.85'
$ws.Range('E2').Value = '  -1.29%  '
$ws.Range('D3').Value = '2.213.08'
$ws.Range('E3').Value = '  -1.49%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '''240.87'
$ws.Range('E5').Value = '  -2.03%  '
$ws.Range('D6').Value = '''0.625'
$ws.Range('E6').Value = '  -0.88%  '
$ws.Range('D7').Value = '''72.64'
$ws.Range('E7').Value = '  -4.14%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').Value = '''0.602'
$ws.Range('E9').Value = '  -2.19%  '
$ws.Range('D10').Value = '''42.18'
$ws.Range('E10').Value = '  -3.49%  '
$ws.Range('D11').Value = '''0.0950'
$ws.Range('E11').Value = '  +0.06%  '
$ws.Range('D12').Value = '''6.96'
$ws.Range('E12').Value = '  -4.63%  '
$ws.Range('E13').Value = '  -0.46%  '
$ws.Range('D14').Value = '2.544.24'
$ws.Range('E14').Value = '  -1.63%  '
$ws.Range('D15').Value = '''14.20'
$ws.Range('E15').Value = '  -2.65%  '
$ws.Range('D16').Value = '''0.833'
$ws.Range('E16').Value = '  -2.46%  '
$ws.Range('D17').Value = '2.203.50'
$ws.Range('E17').Value = '  -1.87%  '
$ws.Range('D18').Value = '41.811.00'
$ws.Range('E18').Value = '  -1.03%  '
$ws.Range('E19').Value = '  +4.48%  '
$ws.Range('D20').Value = '''72.49'
$ws.Range('E20').Value = '  +0.35%  '
$ws.Range('D21').Value = '''6.14'
$ws.Range('E21').Value = '  -0.66%  '
$ws.Range('D22').Value = '''10.68'
$ws.Range('E22').Value = '  +16.88%  '
$ws.Range('D23').Value = '''229.13'
$ws.Range('E23').Value = '  -1.12%  '
$ws.Range('D24').Value = '''2.05'
$ws.Range('E24').Value = '  -7.79%  '
$ws.Range('D25').Value = '''0.998'
$ws.Range('E25').Value = '  -0.24%  '
$ws.Range('D26').Value = '''11.45'
$ws.Range('E26').Value = '  +0.20%  '
$ws.Range('D27').Value = '''3.64'
$ws.Range('E27').Value = '  +0.58%  '
$ws.Range('E28').Value = '  -2.32%  '
$ws.Range('E29').Value = '  -1.02%  '
$ws.Range('D30').Value = '''167.41'
$ws.Range('E30').Value = '  -0.51%  '
$ws.Range('D31').Value = '''20.42'
$ws.Range('E31').Value = '  -1.33%  '
$ws.Range('D32').Value = '''5.60'
$ws.Range('E32').Value = '  +6.34%  '
$ws.Range('D33').Value = '''0.0792'
$ws.Range('E33').Value = '  -4.77%  '
$ws.Range('D34').Value = '''30.26'
$ws.Range('E34').Value = '  -1.84%  '
$ws.Range('E35').Value = '  -1.09%  '
$ws.Range('E36').Value = '  -12.78%  '
$ws.Range('D37').Value = '''4.23'
$ws.Range('E37').Value = '  -6.20%  '
$ws.Range('D38').Value = '''0.0300'
$ws.Range('E38').Value = '  -5.94%  '
$ws.Range('D39').Value = '''13.78'
$ws.Range('E39').Value = '  +1.13%  '
$ws.Range('B40').Value = 'LidoDAOToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D40').Value = '''2.11'
$ws.Range('E40').Value = '  -3.52%  '
$ws.Range('B41').Value = 'MultiversX'
$ws.Range('C41').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D41').Value = '''64.55'
$ws.Range('E41').Value = '  +1.92%  '
$ws.Range('D42').Value = '''5.61'
$ws.Range('E42').Value = '  -3.67%  '
$ws.Range('D43').Value = '''0.196'
$ws.Range('E43').Value = '  -3.33%  '
$ws.Range('D44').Value = '''8.64'
$ws.Range('E44').Value = '  -1.69%  '
$ws.Range('D45').Value = '''103.91'
$ws.Range('E45').Value = '  -3.72%  '
$ws.Range('E46').Value = '  -1.66%  '
$ws.Range('E47').Value = '  -0.88%  '
$ws.Range('D48').Value = '''1.11'
$ws.Range('E48').Value = '  -1.62%  '
$ws.Range('E49').Value = '  -2.25%  '
$ws.Range('E50').Value = '  +0.18%  '
$ws.Range('D51').Value = '2.418.37'
$ws.Range('E51').Value = '  -1.79%  '
